$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace the Kaggle ranking phrase.
#    The exact same phrase appears twice (professional summary paragraph and
#    the "Player Trajectory Prediction" project subtitle) and both need the
#    same replacement text, so a single ReplaceAll covers both occurrences.
# ---------------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.Execute("74th open / 94th closed of 1,134 teams", $true, $false, $false, $false, $false, $true, 1, $false, "Top 8% of 1,134 teams", 2)

# ---------------------------------------------------------------------------
# 2. Update the significance figure in the Missing Persons Outlier Detection
#    bullet (44.75 sigma -> up to 46.86 sigma).
# ---------------------------------------------------------------------------
$sigma = [char]0x03C3
$find2 = $d.Content.Find
$find2.Execute("at 44.75" + $sigma + " significance", $true, $false, $false, $false, $false, $true, 1, $false, "at up to 46.86" + $sigma + " significance", 2)

# ---------------------------------------------------------------------------
# 3. Add the new "OE-OS (In Progress)" project section right after the
#    "Missing Persons Outlier Detection" bullets and before "AI Homelab &
#    Active Memory Network".
# ---------------------------------------------------------------------------
$bullet = [char]0x2022
$nsW = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParaXml = @(
    '<w:p ' + $nsW + '><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">OE-OS (In Progress)</w:t></w:r></w:p>',
    '<w:p ' + $nsW + '><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Distributed AI Orchestration Platform | Python / FastAPI</w:t></w:r></w:p>',
    '<w:p ' + $nsW + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">' + $bullet + ' Three-tier LLM routing (local Ollama to cheap API to premium models) routing ~80% of requests to free local models</w:t></w:r></w:p>',
    '<w:p ' + $nsW + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">' + $bullet + ' Triple-layer RAG memory: BM25 over 5M+ chunks, ChromaDB semantic search, Redis session cache</w:t></w:r></w:p>',
    '<w:p ' + $nsW + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">' + $bullet + ' 18 MCP-compatible tools and multi-agent sandbox with 4 LLM personas at zero API cost</w:t></w:r></w:p>',
    '<w:p ' + $nsW + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="100"/></w:pPr><w:r><w:t xml:space="preserve">' + $bullet + ' 4,200+ lines of async Python on FastAPI for a private multi-node GPU cluster</w:t></w:r></w:p>'
)

# Locate the anchor paragraph ("Built 7-page interactive Streamlit dashboard
# with geospatial visualization") via Find, then use its paragraph Index so
# we don't depend on hard-coded paragraph numbers.
$findAnchor = $d.Content.Find
$findAnchor.Execute("Built 7-page interactive Streamlit dashboard with geospatial visualization", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIndex = $findAnchor.Parent.Paragraphs(1).Index

# Create six new, blank paragraphs right after the anchor paragraph.
$anchor = $d.Paragraphs($anchorIndex)
for ($i = 1; $i -le $newParaXml.Count; $i++) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $d.Paragraphs($anchorIndex + $i)
}

# Fill each new (still-empty) paragraph with its exact target XML so the
# formatting (bold heading, italic subtitle, bulleted list items) matches
# the rest of the resume precisely.
for ($i = 0; $i -lt $newParaXml.Count; $i++) {
    $p = $d.Paragraphs($anchorIndex + 1 + $i)
    $p.Range.InsertXML($newParaXml[$i])
}

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
